$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (matches source formatting,
# which stores prices like "1.002" / "312.72" as literal text, not numbers).
$dCells = @("D2","D3","D4","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D18","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D33","D34","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.233.52"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "1.818.22"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "312.72"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "0.4616"
$ws.Range("E7").Value = "  +5.13%  "
$ws.Range("D8").Value = "0.3751"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").Value = "0.07409"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "0.8677"
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("D11").Value = "20.60"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.823.25"
$ws.Range("E12").Value = "  -4.17%  "
$ws.Range("D13").Value = "6.657"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "5.404"
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("D15").Value = "0.07087"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "92.09"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").Value = "0.000008748"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "27.250.69"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("D22").Value = "5.319"
$ws.Range("E22").Value = "  +3.44%  "
$ws.Range("D23").Value = "10.91"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").Value = "2.049.29"
$ws.Range("E24").Value = "  -5.79%  "
$ws.Range("D25").Value = "1.941"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").Value = "152.26"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").Value = "2.256"
$ws.Range("E27").Value = "  +2.68%  "
$ws.Range("D28").Value = "18.58"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("D29").Value = "5.289"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("D30").Value = "117.08"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "0.08890"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("E32").Value = "  +5.82%  "
$ws.Range("D33").Value = "1.179"
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("D34").Value = "4.523"
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "1.111"
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("D38").Value = "0.01961"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").Value = "0.05247"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("D40").Value = "7.260"
$ws.Range("E40").Value = "  +3.65%  "
$ws.Range("D41").Value = "2.381"
$ws.Range("E41").Value = "  +21.51%  "
$ws.Range("D42").Value = "2.914"
$ws.Range("E42").Value = "  +3.84%  "
$ws.Range("D43").Value = "0.5293"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "0.1687"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "8.590"
$ws.Range("E45").Value = "  +1.82%  "
$ws.Range("D46").Value = "0.5047"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("D47").Value = "10.45"
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").Value = "105.08"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "1.672"
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "0.06328"
$ws.Range("E51").Value = "  +0.31%  "
